# Cambiamos el LDR por un sensor infrarojo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "LDR" component is replaced by an infrared sensor
$ws.Range("B5").Value = " INFRAROGO"

# Its price is no longer computed from quantity (25*C5); it's now a flat value
$ws.Range("D5").Value = 400

# A new (blank-looking) row label appears right above the "Total" row
$ws.Range("B12").Value = "  "

# A new, very wide column E is introduced
$ws.Columns("E").ColumnWidth = 131.62

# Recalculate dependent totals (D13 = SUM(D5:D12), D14 = D13/4)
$wb.Application.Calculate()

# Leave the active selection on E3, matching the edited workbook
$ws.Range("E3").Select()
